# "add total yield col. info"
#
# Populates the "eultion_vol_ul" (L) and "total-yield_ng" (M) columns
# on Sheet1 for the rows that belong to qubit runs whose elution volume
# is now known, adding M = (F)*(L-G)  [ (tube conc) * (elution vol - dilution factor) ].
# A handful of rows in the 255-288 block are also missing a dilution
# factor (G) of 1, which is backfilled at the same time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Fill-YieldRows($FirstRow, $LastRow, $ElutionVol, $NeedsDilutionFactor) {
    for ($r = $FirstRow; $r -le $LastRow; $r++) {
        if ($NeedsDilutionFactor) {
            $ws.Cells.Item($r, 7).Value = 1   # column G: dilution_factor
        }
        $ws.Cells.Item($r, 12).Value = $ElutionVol   # column L: eultion_vol_ul
        $ws.Range("M$r").Formula = "=(F$r)*(L$r-G$r)"  # column M: total-yield_ng
    }
}

# Rows 2-11: elution volume 50 ul, dilution factor already present.
Fill-YieldRows 2 11 50 $false

# Rows 170-209: elution volume 14 ul, dilution factor already present.
Fill-YieldRows 170 209 14 $false

# Rows 234-235 and 238-245: elution volume 14 ul, dilution factor already present
# (rows 236-237 are a different record shape and are intentionally skipped).
Fill-YieldRows 234 235 14 $false
Fill-YieldRows 238 245 14 $false

# Rows 255-288: elution volume 14 ul; these rows never had a dilution
# factor recorded, so backfill G=1 as well.
Fill-YieldRows 255 288 14 $true

# Rows 289-300: elution volume 14 ul, dilution factor already present
# (L/M cells already existed, just unfilled).
Fill-YieldRows 289 300 14 $false

# Move the active selection the way the author left it.
$null = $ws.Range("J301").Select()
